# Apply the edits described by the commit "fixed Excel resource file" to
# Resources/Risky_Business_Tables.xlsx.
#
# Changes made:
#  1. Add a new total row: C25 holds a raw value (17104) and C24 holds a
#     formula that divides a constant by it (=6789/C25).
#  2. Highlight four cells (B4, C5, C8, C9) with a solid light-green fill
#     (RGB 146, 208, 80 / hex 92D050) - this introduces a new fill + cell
#     style in the workbook.
#  3. Remove the color-scale conditional formatting that used to sit over
#     the B2:E2, B6:E6 and B10:E10 ranges (the rules over B14:C14 and below
#     are left untouched).
#  4. Leave the selection on the newly added C25 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New total row / formula -------------------------------------------
$ws.Range("C25").Value = 17104
$ws.Range("C24").Formula = "=6789/C25"

# --- 2. Highlight cells with a solid green fill ----------------------------
# Interior.Color takes a BGR-packed long; 146,208,80 (0x92,0xD0,0x50) packs
# to 0x50D092 = 5296274.
$green = 5296274
$ws.Range("B4").Interior.Color = $green
$ws.Range("C5").Interior.Color = $green
$ws.Range("C8").Interior.Color = $green
$ws.Range("C9").Interior.Color = $green

# --- 3. Drop the conditional formatting on the three score rows -----------
$ws.Range("B2:E2").FormatConditions.Delete()
$ws.Range("B6:E6").FormatConditions.Delete()
$ws.Range("B10:E10").FormatConditions.Delete()

# --- 4. Move the active selection to the new cell --------------------------
$ws.Range("C25").Select()
